# Update gh-pages to output generated at 456a3b4
# Applies refreshed "想去人数" (F) / "最低票价" (G) values across the
# 展览, 演出 and 全部类型 sheets.

$wb = $excel.ActiveWorkbook

# --- 展览 (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6804
$ws1.Range("F3").Value = 0
$ws1.Range("F5").Value = 439
$ws1.Range("F6").Value = 0
$ws1.Range("F7").Value = 6438
$ws1.Range("G7").Value = 68
$ws1.Range("F8").Value = 56
$ws1.Range("F9").Value = 0
$ws1.Range("F10").Value = 1278
$ws1.Range("F12").Value = 105
$ws1.Range("F13").Value = 0
$ws1.Range("F16").Value = 0
$ws1.Range("F18").Value = 0
$ws1.Range("F19").Value = 4779
$ws1.Range("F20").Value = 81
$ws1.Range("F21").Value = 61
$ws1.Range("F22").Value = 0
$ws1.Range("F24").Value = 0

# --- 演出 (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 0

# --- 全部类型 (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6804
$ws4.Range("F5").Value = 439
$ws4.Range("F6").Value = 146
$ws4.Range("F7").Value = 6438
$ws4.Range("G7").Value = 68
$ws4.Range("F8").Value = 56
$ws4.Range("F9").Value = 0
$ws4.Range("F10").Value = 1278
$ws4.Range("F11").Value = 0
$ws4.Range("F12").Value = 105
$ws4.Range("F13").Value = 394
$ws4.Range("F15").Value = 18
$ws4.Range("F16").Value = 0
$ws4.Range("F17").Value = 45
$ws4.Range("F18").Value = 0
$ws4.Range("F19").Value = 0
$ws4.Range("F20").Value = 0
$ws4.Range("F21").Value = 0
$ws4.Range("F22").Value = 61
$ws4.Range("F23").Value = 248
$ws4.Range("F25").Value = 129
